$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Title ---
Replace-Text "Human Space Exploration: Endeavors and Prospects" "Exploring the Fascinating World of Biology: Unveiling the Secrets of Life"

# --- Author name ---
Replace-Text " Samuel Harrison" " Erika Johnson"

# --- Author email (three runs: "harrison" / "." / "samuel@aereospace" / "." / "com") ---
Replace-Text "harrison.samuel@aereospace.com" "erika.johnson@academicschool.edu"

# --- Main body paragraph (three sentences-groups separated by manual line breaks) ---
$vv = "`v`v"
$bodyOld = "As humanity seeks to broaden its horizons beyond the confines of our earthly home, the realm of space exploration stands as a testament to our insatiable curiosity, ingenuity, and determination." +
           " With each mission embarked upon, we reaffirm our commitment to unraveling the mysteries of the universe and unlocking its boundless potential." +
           " From celestial navigation by ancient civilizations to the ambitious voyages of the modern era, space exploration has ignited our imaginations and shaped our understanding of the cosmos." +
           $vv +
           "Driven by an insatiable thirst for knowledge, we strive to decipher the origins of life, the formation of planets, and the fundamental forces that govern reality." +
           " Our expeditions into space have yielded invaluable insights into Earth's intricate climate systems, aided in predicting natural disasters, and enhanced communication networks." +
           " We have unearthed resources and minerals with the potential to revolutionize industries, leading to advancements in technology, medicine, and renewable energy sources." +
           $vv +
           "Our space endeavors have also kindled collaborations among nations, transcending political and cultural boundaries." +
           " Boldly embarking on shared missions, countries have forged alliances, fostering peace and cooperation on Earth while collectively venturing towards the cosmos." +
           " These collaborative efforts have deepened our collective understanding and ignited a sense of global unity, highlighting the transformative power of exploration."

$bodyNew = "Biology, the study of life and its processes, is an intriguing and dynamic field that unlocks the mysteries of the living world." +
           " From the smallest microorganisms to the grandest ecosystems, biology unveils the intricate patterns and principles that govern the functioning and diversity of all living things." +
           $vv +
           "In this captivating journey into the depths of biology, we will delve into the fundamental units of life: cells." +
           " These microscopic entities, despite their diminutive size, are marvels of complexity, carrying out intricate processes that support life." +
           " We will discover the remarkable diversity of organisms, ranging from single-celled bacteria to multicellular organisms with astonishing adaptations, and explore the delicate balance of ecosystems that sustain life on Earth." +
           $vv +
           "As we unravel the intricacies of biological systems, we will uncover the significance of DNA and genetic information in shaping the characteristics and traits of organisms." +
           " We will unravel the fascinating mechanisms of evolution, the driving force behind the diversity of life, and delve into the extraordinary adaptations that enable organisms to thrive in a multitude of environments."

Replace-Text $bodyOld $bodyNew

# --- Summary paragraph ---
$summaryOld = "Human space exploration is a captivating tapestry woven with threads of curiosity, ambition, and human ingenuity." +
              " Our forays into space have yielded critical scientific knowledge, technological advancements, and resources that have indelibly impacted life on Earth." +
              " As we look towards future missions, we are filled with hope and determination to push the boundaries of our understanding, unravel the mysteries of the universe, and secure the future of our species among the myriad celestial bodies that dance in the blackness of space."

$summaryNew = "Biology, the study of life and its processes, offers a profound understanding of the marvelous diversity, intricacies, and interconnectedness of living organisms." +
              " It encompasses the exploration of cells, the fundamental units of life, and unravels the significance of DNA and genetic information in shaping the characteristics and traits of organisms." +
              " Biology unveils the mechanisms of evolution, the driving force behind the diversity of life." +
              " It delves into the remarkable adaptations that enable organisms to thrive in various environments and explores the delicate balance of ecosystems that sustain life on Earth, highlighting the interconnectedness of all living things." +
              " Biology provides a profound appreciation for the beauty and complexity of life, inspiring awe and a sense of stewardship for the natural world." +
              "`r"

Replace-Text $summaryOld $summaryNew

Write-Output "done"
